$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'316.87"
$ws.Range("E2").Value = "'2.17%"
$ws.Range("D3").Value = "'41.10"
$ws.Range("E3").Value = "'0.05%"
$ws.Range("D4").Value = "'5.166"
$ws.Range("E4").Value = "'0.80%"
$ws.Range("D5").Value = "'0.07635"
$ws.Range("E5").Value = "'-0.71%"
$ws.Range("B6").Value = "FTXToken"
$ws.Range("C6").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D6").Value = "'1.683"
$ws.Range("E6").Value = "'3.68%"
$ws.Range("B7").Value = "MXToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D7").Value = "'0.9324"
$ws.Range("E7").Value = "'1.33%"
$ws.Range("B8").Value = "BTSEToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D8").Value = "'2.425"
$ws.Range("E8").Value = "'-1.62%"
$ws.Range("B9").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C9").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D9").Value = "'0.1248"
$ws.Range("E9").Value = "'2.39%"
$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D10").Value = "'0.1825"
$ws.Range("E10").Value = "'-0.60%"
$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D11").Value = "'0.09034"
$ws.Range("E11").Value = "'-0.74%"
$ws.Range("B12").Value = "BitrueCoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D12").Value = "'0.04164"
$ws.Range("E12").Value = "'-2.19%"
$ws.Range("B13").Value = "BitMartToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D13").Value = "'0.1056"
$ws.Range("E13").Value = "'0.56%"
$ws.Range("B14").Value = "BitForexToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D14").Value = "'0.001268"
$ws.Range("E14").Value = "'1.94%"
$ws.Range("B15").Value = "TigerCash"
$ws.Range("C15").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D15").Value = "'0.005918"
$ws.Range("E15").Value = "'1.55%"
$ws.Range("B16").Value = "UpBots"
$ws.Range("C16").Value = "https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt"
$ws.Range("D16").Value = "'0.007491"
$ws.Range("E16").Value = "'1,897.31%"
$ws.Range("B17").Value = "LEO"
$ws.Range("C17").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D17").Value = "'3.351"
$ws.Range("E17").Value = "'-0.06%"
$ws.Range("B18").Value = "GateToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D18").Value = "'4.315"
$ws.Range("E18").Value = "'0.91%"
$ws.Range("D19").Value = "'0.3361"
$ws.Range("E19").Value = "'1.45%"
$ws.Range("D20").Value = "'8.359"
$ws.Range("E20").Value = "'21.04%"
$ws.Range("D21").Value = "'0.1347"
$ws.Range("E21").Value = "'-2.31%"
$ws.Range("E22").Value = "'7.46%"
$ws.Range("E23").Value = "'0.05%"
$ws.Range("D24").Value = "'0.001272"
$ws.Range("E24").Value = "'0.88%"
$ws.Range("D25").Value = "'0.004083"
$ws.Range("E25").Value = "'0.17%"
$ws.Range("D26").Value = "'0.0001276"
$ws.Range("E26").Value = "'0.67%"
$ws.Range("D38").Value = "'0.02484"
$ws.Range("E38").Value = "'0.42%"
$ws.Range("D39").Value = "'0.05253"
$ws.Range("E39").Value = "'-0.27%"
$ws.Range("D40").Value = "'0.007786"
$ws.Range("E40").Value = "'-0.53%"
$ws.Range("D41").Value = "'0.1298"
$ws.Range("E41").Value = "'-1.16%"
$ws.Range("D42").Value = "'0.007068"
$ws.Range("E42").Value = "'4.08%"
$ws.Range("D43").Value = "'0.002098"
$ws.Range("E43").Value = "'14.11%"
$ws.Range("D44").Value = "'0.008239"
$ws.Range("D45").Value = "'0.3430"
$ws.Range("E45").Value = "'10.74%"
$ws.Range("D46").Value = "'0.00006685"
$ws.Range("E46").Value = "'-2.23%"
$ws.Range("D47").Value = "'0.00000000753"
$ws.Range("E47").Value = "'0.66%"
$ws.Range("D48").Value = "'0.2189"
$ws.Range("E48").Value = "'-4.25%"
$ws.Range("D49").Value = "'0.004219"
$ws.Range("E49").Value = "'3.12%"
$ws.Range("D50").Value = "'0.00002109"
$ws.Range("E50").Value = "'0.66%"
$ws.Range("D51").Value = "'0.0002009"
$ws.Range("E51").Value = "'0.66%"
